$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 4
$ws.Range("H4").Value = 176.25
$ws.Range("I4").Value = 74.5
$ws.Range("J4").Value = 278
$ws.Range("K4").Value = 74.5
$ws.Range("L4").Value = 278
$ws.Range("M4").Value = 39.5
$ws.Range("N4").Value = -506

# ALC row 7
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

# ALC row 14
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

# ALC row 16
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 1000
$ws.Range("K16").Value = 1000
$ws.Range("M16").Value = -770

# ALC row 54
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

# ALC row 70
$ws.Range("H70").Value = 5866
$ws.Range("I70").Value = 3800
$ws.Range("J70").Value = 9998
$ws.Range("K70").Value = 11400
$ws.Range("L70").Value = 29994
$ws.Range("M70").Value = -11130
$ws.Range("N70").Value = -30534

# ALC row 73
$ws.Range("H73").Value = 5866
$ws.Range("I73").Value = 3800
$ws.Range("J73").Value = 9998
$ws.Range("K73").Value = 11400
$ws.Range("L73").Value = 29994
$ws.Range("M73").Value = -10464
$ws.Range("N73").Value = -31866

# ALC row 74
$ws.Range("H74").Value = 4500
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 5000
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -4064
$ws.Range("N74").Value = -5872

# ALC row 77
$ws.Range("H77").Value = 4500
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 25000
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -20320
$ws.Range("N77").Value = -29360

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws.Range("H2").Value = 5000
$ws.Range("J2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("N2").Value = -5226

# ARM row 5
$ws.Range("H5").Value = 95.7
$ws.Range("I5").Value = 88.14286
$ws.Range("K5").Value = 88.14286
$ws.Range("M5").Value = 23.85714

# ARM row 6
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

# ARM row 45
$ws.Range("H45").Value = 3508.5
$ws.Range("I45").Value = 3371.4443
$ws.Range("K45").Value = 3371.4443
$ws.Range("M45").Value = -2994.4443

# ARM row 97
$ws.Range("H97").Value = 605
$ws.Range("I97").Value = 605
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 605
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -109
$ws.Range("N97").ClearContents()

# ARM row 116
$ws.Range("H116").Value = 5000
$ws.Range("J116").Value = 5000
$ws.Range("L116").Value = 5000
$ws.Range("N116").Value = -9588

# ARM row 132
$ws.Range("H132").Value = 1139.3334
$ws.Range("I132").Value = 1139.3334
$ws.Range("K132").Value = 3418.0002
$ws.Range("M132").Value = -888.0001999999999

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws.Range("H3").Value = 5000
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("N3").Value = -5228

# BSM row 4
$ws.Range("H4").Value = 95.7
$ws.Range("I4").Value = 88.14286
$ws.Range("K4").Value = 88.14286
$ws.Range("M4").Value = 26.85714

# BSM row 94
$ws.Range("H94").Value = 2100
$ws.Range("I94").Value = 2100
$ws.Range("K94").Value = 2100
$ws.Range("M94").Value = -1649

# BSM row 105
$ws.Range("H105").Value = 1945.75
$ws.Range("I105").Value = 1945.75
$ws.Range("K105").Value = 1945.75
$ws.Range("M105").Value = -198.75

# BSM row 134
$ws.Range("H134").Value = 2022.5
$ws.Range("I134").Value = 1458.3334
$ws.Range("J134").Value = 3715
$ws.Range("K134").Value = 4375.0002
$ws.Range("L134").Value = 11145
$ws.Range("M134").Value = -1840.0002
$ws.Range("N134").Value = -16215

$ws = $wb.Worksheets.Item("CRP")
# CRP row 23
$ws.Range("H23").Value = 1009
$ws.Range("I23").Value = 1009
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 1009
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -769
$ws.Range("N23").ClearContents()

# CRP row 27
$ws.Range("H27").Value = 1009
$ws.Range("I27").Value = 1009
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1009
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -817
$ws.Range("N27").ClearContents()

# CRP row 96
$ws.Range("H96").Value = 20000
$ws.Range("J96").Value = 20000
$ws.Range("L96").Value = 20000
$ws.Range("N96").Value = -25492

# CRP row 119
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# CUL row 4
$ws.Range("H4").Value = 250449.75
$ws.Range("I4").Value = 500499.5
$ws.Range("K4").Value = 1501498.5
$ws.Range("M4").Value = -1501386.5

# CUL row 75
$ws.Range("H75").Value = 1012.5
$ws.Range("I75").Value = 25
$ws.Range("J75").Value = 2000
$ws.Range("K75").Value = 75
$ws.Range("L75").Value = 6000
$ws.Range("M75").Value = 923
$ws.Range("N75").Value = -7996

# CUL row 78
$ws.Range("H78").Value = 1012.5
$ws.Range("I78").Value = 25
$ws.Range("J78").Value = 2000
$ws.Range("K78").Value = 225
$ws.Range("L78").Value = 18000
$ws.Range("M78").Value = 4767
$ws.Range("N78").Value = -27984

$ws = $wb.Worksheets.Item("GSM")
# GSM row 5
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

# GSM row 53
$ws.Range("H53").Value = 21916.334
$ws.Range("J53").Value = 22874.5
$ws.Range("L53").Value = 22874.5
$ws.Range("N53").Value = -24136.5

$ws = $wb.Worksheets.Item("LTW")
# LTW row 46
$ws.Range("H46").Value = 4372.1816
$ws.Range("I46").Value = 4480.75
$ws.Range("J46").Value = 4082.6667
$ws.Range("K46").Value = 4480.75
$ws.Range("L46").Value = 4082.6667
$ws.Range("M46").Value = -4292.75
$ws.Range("N46").Value = -4458.6667

# LTW row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

# LTW row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# WVR row 109
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()
